# Apply the "dynamically update the punch inner and outer section" edit:
# new mesh/input-file rows describing the split punch boundary, and the
# Meshes sheet becomes the active sheet/tab.

$wb = $excel.ActiveWorkbook

$description = "the punch" + [char]0x2019 + "s outer boundary was split to allow for dirichlet temperature boundary condition to be applied to upper part and thermal contact condition to be applied on the lower part which is inside the die."

# --- Input-files sheet -----------------------------------------------
$wsInput = $wb.Worksheets.Item("Input-files")

$wsInput.Range("A5").Value = "pellet_v5_20180527.i"
$wsInput.Range("B5").Value = $description

$wsInput.Range("A7").Value = "pellet_v7_20180706.i"

$wsInput.Activate()
$wsInput.Range("A5").Select() | Out-Null

# --- Meshes sheet ------------------------------------------------------
$wsMeshes = $wb.Worksheets.Item("Meshes")

$wsMeshes.Range("A4").Value = "pellet_v4_20180517.e"

$wsMeshes.Range("A5").Value = "pellet_v5_20180529.e"
$wsMeshes.Range("B5").Value = $description

$wsMeshes.Activate()
$wsMeshes.Range("B9").Select() | Out-Null
